$wb = $excel.ActiveWorkbook

# Remove tabSelected from the current active sheet (clear_data) by adding new sheets after it,
# which will automatically become the new active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$studentBook = $wb.Worksheets.Add($null, $lastSheet)
$studentBook.Name = "student_book"

$studentBook.Range('A1').Value = 'วิชาคณิตศาสตร์'
$studentBook.Range('A2').Value = 'r_m_1_1_1'
$studentBook.Range('B2').Value = 'ชั้นประถมศีกษาปีที่ 1'
$studentBook.Range('C2').Value = 'หนังสือเรียน'
$studentBook.Range('A3').Value = 'r_m_1_2_1'
$studentBook.Range('C3').Value = 'แบบฝึกทักษะ เตรียมความพร้อม'
$studentBook.Range('A4').Value = 'r_m_1_3_1'
$studentBook.Range('C4').Value = 'แบบฝึกทักษะ เล่ม 1'
$studentBook.Range('A5').Value = 'r_m_1_4_1'
$studentBook.Range('C5').Value = 'แบบฝึกทักษะ เล่ม 2'
$studentBook.Range('A6').Value = 'r_m_2_1_1'
$studentBook.Range('B6').Value = 'ชั้นประถมศีกษาปีที่ 2'
$studentBook.Range('C6').Value = 'หนังสือเรียน'
$studentBook.Range('A7').Value = 'r_m_2_2_1'
$studentBook.Range('C7').Value = 'แบบฝึกทักษะ เล่ม 1'
$studentBook.Range('A8').Value = 'r_m_2_3_1'
$studentBook.Range('C8').Value = 'แบบฝึกทักษะ เล่ม 2'
$studentBook.Range('A9').Value = 'r_m_3_1_1'
$studentBook.Range('B9').Value = 'ชั้นประถมศีกษาปีที่ 3'
$studentBook.Range('C9').Value = 'หนังสือเรียน'
$studentBook.Range('A10').Value = 'r_m_3_2_1'
$studentBook.Range('C10').Value = 'แบบฝึกทักษะ เล่ม 1'
$studentBook.Range('A11').Value = 'r_m_3_3_1'
$studentBook.Range('C11').Value = 'แบบฝึกทักษะ เล่ม 2'
$studentBook.Range('A12').Value = 'r_m_4_1_1'
$studentBook.Range('B12').Value = 'ชั้นประถมศีกษาปีที่ 4'
$studentBook.Range('C12').Value = 'หนังสือเรียน'
$studentBook.Range('A13').Value = 'r_m_4_2_1'
$studentBook.Range('C13').Value = 'แบบฝึกทักษะ เล่ม 1'
$studentBook.Range('A14').Value = 'r_m_4_3_1'
$studentBook.Range('C14').Value = 'แบบฝึกทักษะ เล่ม 2'
$studentBook.Range('A15').Value = 'r_m_5_1_1'
$studentBook.Range('B15').Value = 'ชั้นประถมศีกษาปีที่ 5'
$studentBook.Range('C15').Value = 'หนังสือเรียน'
$studentBook.Range('A16').Value = 'r_m_5_2_1'
$studentBook.Range('C16').Value = 'แบบฝึกทักษะ เล่ม 1'
$studentBook.Range('A17').Value = 'r_m_5_3_1'
$studentBook.Range('C17').Value = 'แบบฝึกทักษะ เล่ม 2'
$studentBook.Range('A18').Value = 'r_m_6_1_1'
$studentBook.Range('B18').Value = 'ชั้นประถมศีกษาปีที่ 6'
$studentBook.Range('C18').Value = 'หนังสือเรียน'
$studentBook.Range('A19').Value = 'r_m_6_2_1'
$studentBook.Range('C19').Value = 'แบบฝึกทักษะ เล่ม 1'
$studentBook.Range('A20').Value = 'r_m_6_3_1'
$studentBook.Range('C20').Value = 'แบบฝึกทักษะ เล่ม 2'
$studentBook.Range('A21').Value = 'r_m_7_1_1'
$studentBook.Range('B21').Value = 'ชั้นมัธยมศึกษาปีที่ 1'
$studentBook.Range('C21').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 1'
$studentBook.Range('A22').Value = 'r_m_7_2_1'
$studentBook.Range('C22').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 2'
$studentBook.Range('A23').Value = 'r_m_7_3_1'
$studentBook.Range('C23').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 1'
$studentBook.Range('A24').Value = 'r_m_7_4_1'
$studentBook.Range('C24').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 2'
$studentBook.Range('A25').Value = 'r_m_8_1_1'
$studentBook.Range('B25').Value = 'ชั้นมัธยมศึกษาปีที่ 2'
$studentBook.Range('C25').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 1'
$studentBook.Range('A26').Value = 'r_m_8_2_1'
$studentBook.Range('C26').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 2'
$studentBook.Range('A27').Value = 'r_m_8_3_1'
$studentBook.Range('C27').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 1'
$studentBook.Range('A28').Value = 'r_m_8_4_1'
$studentBook.Range('C28').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 2'
$studentBook.Range('A29').Value = 'r_m_9_1_1'
$studentBook.Range('B29').Value = 'ชั้นมัธยมศึกษาปีที่ 3'
$studentBook.Range('C29').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 1'
$studentBook.Range('A30').Value = 'r_m_9_2_1'
$studentBook.Range('C30').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 2'
$studentBook.Range('A31').Value = 'r_m_9_3_1'
$studentBook.Range('C31').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 1'
$studentBook.Range('A32').Value = 'r_m_9_4_1'
$studentBook.Range('C32').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 2'
$studentBook.Range('A33').Value = 'r_m_101112_1_1'
$studentBook.Range('B33').Value = 'ชั้นมัธยมศึกษาปีที่ 4-6'
$studentBook.Range('C33').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 1'
$studentBook.Range('A34').Value = 'r_m_101112_2_1'
$studentBook.Range('C34').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 2'
$studentBook.Range('A35').Value = 'r_m_101112_3_1'
$studentBook.Range('C35').Value = 'หนังสือเรียน รายวิชาพื้นฐาน เล่ม 3'
$studentBook.Range('A36').Value = 'r_m_101112_4_1'
$studentBook.Range('C36').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 1'
$studentBook.Range('A37').Value = 'r_m_101112_5_1'
$studentBook.Range('C37').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 2'
$studentBook.Range('A38').Value = 'r_m_101112_6_1'
$studentBook.Range('C38').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 3'
$studentBook.Range('A39').Value = 'r_m_101112_7_1'
$studentBook.Range('C39').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 4'
$studentBook.Range('A40').Value = 'r_m_101112_8_1'
$studentBook.Range('C40').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 5'
$studentBook.Range('A41').Value = 'r_m_101112_9_1'
$studentBook.Range('C41').Value = 'หนังสือเรียน รายวิชาเพิ่มเติม เล่ม 6'

$studentBook.Columns.Item(1).ColumnWidth = 30.833333333333332
$studentBook.Columns.Item(2).ColumnWidth = 52
$studentBook.Columns.Item(3).ColumnWidth = 44
$studentBook.Range("D6").Select()

$teacherBook = $wb.Worksheets.Add($null, $studentBook)
$teacherBook.Name = "teacher_book"

$teacherBook.Range('A1').Value = 'code'
$teacherBook.Range('B1').Value = 'grade'
$teacherBook.Range('C1').Value = 'book'
$teacherBook.Range('A2').Value = 'คู่มือครูวิชาวิทยาศาสตร์'
$teacherBook.Range('A3').Value = 'r_s_ins_1_1_1'
$teacherBook.Range('B3').Value = 'ชั้นประถมศีกษาปีที่ 1'
$teacherBook.Range('C3').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('A4').Value = 'r_s_ins_2_1_1'
$teacherBook.Range('B4').Value = 'ชั้นประถมศีกษาปีที่ 2'
$teacherBook.Range('C4').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('A5').Value = 'r_s_ins_3_1_1'
$teacherBook.Range('B5').Value = 'ชั้นประถมศีกษาปีที่ 3'
$teacherBook.Range('C5').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('A6').Value = 'r_s_ins_4_1_1'
$teacherBook.Range('B6').Value = 'ชั้นประถมศีกษาปีที่ 4'
$teacherBook.Range('C6').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('A7').Value = 'r_s_ins_5_1_1'
$teacherBook.Range('B7').Value = 'ชั้นประถมศีกษาปีที่ 5'
$teacherBook.Range('C7').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('A8').Value = 'r_s_ins_6_1_1'
$teacherBook.Range('B8').Value = 'ชั้นประถมศีกษาปีที่ 6'
$teacherBook.Range('C8').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('A9').Value = 'r_s_ins_7_1_1'
$teacherBook.Range('B9').Value = 'ชั้นมัธยมศึกษาปีที่ 1'
$teacherBook.Range('C9').Value = 'คู่มือครู รายวิชาพื้นฐาน เล่ม 1'
$teacherBook.Range('A10').Value = 'r_s_ins_7_2_1'
$teacherBook.Range('C10').Value = 'คู่มือครู รายวิชาพื้นฐาน เล่ม 2'
$teacherBook.Range('A11').Value = 'r_s_ins_8_1_1'
$teacherBook.Range('B11').Value = 'ชั้นมัธยมศึกษาปีที่ 2'
$teacherBook.Range('C11').Value = 'คู่มือครู รายวิชาพื้นฐาน เล่ม 1'
$teacherBook.Range('A12').Value = 'r_s_ins_8_2_1'
$teacherBook.Range('C12').Value = 'คู่มือครู รายวิชาพื้นฐาน เล่ม 2'
$teacherBook.Range('A13').Value = 'r_s_ins_9_1_1'
$teacherBook.Range('B13').Value = 'ชั้นมัธยมศึกษาปีที่ 3'
$teacherBook.Range('C13').Value = 'คู่มือครู รายวิชาพื้นฐาน เล่ม 1'
$teacherBook.Range('A14').Value = 'r_s_ins_9_2_1'
$teacherBook.Range('C14').Value = 'คู่มือครู รายวิชาพื้นฐาน เล่ม 2'
$teacherBook.Range('A15').Value = 'r_s_ins_789_additional_1_1'
$teacherBook.Range('B15').Value = 'ชั้นมัธยมศึกษาปีที่ 1-3'
$teacherBook.Range('C15').Value = 'คู่มือครู เชื้อเพลิงเพื่อการคมนาคม'
$teacherBook.Range('A16').Value = 'r_s_ins_789_additional_2_1'
$teacherBook.Range('C16').Value = 'คู่มือครู ของเล่นเชิงวิทยาศาสตร์'
$teacherBook.Range('A17').Value = 'r_s_ins_789_additional_3_1'
$teacherBook.Range('C17').Value = 'คู่มือครู วิทยาศาสตร์กับความงาม'
$teacherBook.Range('A18').Value = 'r_s_ins_789_additional_4_1'
$teacherBook.Range('C18').Value = 'คู่มือครู สนุกกับโครงงานวิทยาศาสตร์'
$teacherBook.Range('A19').Value = 'r_s_ins_789_additional_5_1'
$teacherBook.Range('C19').Value = 'คู่มือครู พลังงานทดแทนกับการใช้ประโยชน์'
$teacherBook.Range('A20').Value = 'r_s_ins_101112n_1_1'
$teacherBook.Range('B20').Value = 'ชั้นมัธยมศึกษาปีที่ 4-6'
$teacherBook.Range('C20').Value = 'คู่มือครู การเคลื่อนที่และแรงในธรรมชาติ'
$teacherBook.Range('D20').Value = 'วิทยาศาสตร์ สำหรับนักเรียนที่ไม่เน้นวิทยาศาสตร์ '
$teacherBook.Range('A21').Value = 'r_s_ins_101112n_2_1'
$teacherBook.Range('C21').Value = 'คู่มือครู ดวงดาวและโลกของเรา'
$teacherBook.Range('A22').Value = 'r_s_ins_101112n_3_1'
$teacherBook.Range('C22').Value = 'คู่มือครู พลังงาน'
$teacherBook.Range('A23').Value = 'r_s_ins_101112n_4_1'
$teacherBook.Range('C23').Value = 'คู่มือครู พันธุกรรมและสิ่งแวดล้อม'
$teacherBook.Range('A24').Value = 'r_s_ins_101112n_5_1'
$teacherBook.Range('C24').Value = 'คู่มือครู สารและสมบัติของสาร'
$teacherBook.Range('A25').Value = 'r_s_ins_101112n_6_1'
$teacherBook.Range('C25').Value = 'คู่มือครู ดุลยภาพของสิ่งมีชีวิต'
$teacherBook.Range('A26').Value = 'r_s_ins_101112p_1_1'
$teacherBook.Range('C26').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('D26').Value = 'ฟิสิกส์ สำหรับนักเรียนที่เน้นวิทยาศาสตร์'
$teacherBook.Range('A27').Value = 'r_s_ins_101112p_2_1'
$teacherBook.Range('C27').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 1'
$teacherBook.Range('A28').Value = 'r_s_ins_101112p_3_1'
$teacherBook.Range('C28').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 2'
$teacherBook.Range('A29').Value = 'r_s_ins_101112p_4_1'
$teacherBook.Range('C29').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 3 '
$teacherBook.Range('A30').Value = 'r_s_ins_101112p_5_1'
$teacherBook.Range('C30').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 4'
$teacherBook.Range('A31').Value = 'r_s_ins_101112p_6_1'
$teacherBook.Range('C31').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 5'
$teacherBook.Range('A32').Value = 'r_s_ins_101112c_1_1'
$teacherBook.Range('C32').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('D32').Value = 'เคมี สำหรับนักเรียนที่เน้นวิทยาศาสตร์'
$teacherBook.Range('A33').Value = 'r_s_ins_101112c_2_1'
$teacherBook.Range('C33').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 1'
$teacherBook.Range('A34').Value = 'r_s_ins_101112c_3_1'
$teacherBook.Range('C34').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 2'
$teacherBook.Range('A35').Value = 'r_s_ins_101112c_4_1'
$teacherBook.Range('C35').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 3'
$teacherBook.Range('A36').Value = 'r_s_ins_101112c_5_1'
$teacherBook.Range('C36').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 4 '
$teacherBook.Range('A37').Value = 'r_s_ins_101112c_6_1'
$teacherBook.Range('C37').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 5 '
$teacherBook.Range('A38').Value = 'r_s_ins_101112b_1_1'
$teacherBook.Range('C38').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('D38').Value = 'ชีววิทยา สำหรับนักเรียนที่เน้นวิทยาศาสตร์'
$teacherBook.Range('A39').Value = 'r_s_ins_101112b_2_1'
$teacherBook.Range('C39').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 1'
$teacherBook.Range('A40').Value = 'r_s_ins_101112b_3_1'
$teacherBook.Range('C40').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 2'
$teacherBook.Range('A41').Value = 'r_s_ins_101112b_4_1'
$teacherBook.Range('C41').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 3'
$teacherBook.Range('A42').Value = 'r_s_ins_101112b_5_1'
$teacherBook.Range('C42').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 4'
$teacherBook.Range('A43').Value = 'r_s_ins_101112b_6_1'
$teacherBook.Range('C43').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 5 '
$teacherBook.Range('A44').Value = 'r_s_ins_101112e_1_1'
$teacherBook.Range('C44').Value = 'คู่มือครู รายวิชาพื้นฐาน'
$teacherBook.Range('D44').Value = 'โลก ดาราศาสตร์ และอวกาศ สำหรับนักเรียนที่เน้นวิทยาศาสตร์'
$teacherBook.Range('A45').Value = 'r_s_ins_101112e_2_1'
$teacherBook.Range('C45').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 1'
$teacherBook.Range('A46').Value = 'r_s_ins_101112e_3_1'
$teacherBook.Range('C46').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 2'
$teacherBook.Range('A47').Value = 'r_s_ins_101112e_4_1'
$teacherBook.Range('C47').Value = 'คู่มือครู รายวิชาเพิ่มเติม เล่ม 3'

$teacherBook.Columns.Item(1).ColumnWidth = 30.833333333333332
$teacherBook.Columns.Item(2).ColumnWidth = 33
$teacherBook.Columns.Item(3).ColumnWidth = 53
$teacherBook.Columns.Item(4).ColumnWidth = 51
$teacherBook.Range("D11").Select()